# Update the "想去人数" (want-to-go count) figures in the 广州-漫展信息 workbook.
# These are the same underlying events re-scraped at a later time, so the
# counts are bumped (slightly) across three sheets: 展览, 演出 and 全部类型.
# (本地生活 is untouched.)

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 870
$ws.Range("F3").Value = 1445
$ws.Range("F4").Value = 1097
$ws.Range("F5").Value = 521
$ws.Range("F6").Value = 222
$ws.Range("F8").Value = 673
$ws.Range("F9").Value = 250
$ws.Range("F13").Value = 151
$ws.Range("F14").Value = 1990
$ws.Range("F15").Value = 430
$ws.Range("F17").Value = 496
$ws.Range("F18").Value = 270
$ws.Range("F22").Value = 663
$ws.Range("F23").Value = 51
$ws.Range("F24").Value = 241
$ws.Range("F25").Value = 960
$ws.Range("F27").Value = 1561
$ws.Range("F28").Value = 303

# --- 演出 (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 282

# --- 全部类型 (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 870
$ws.Range("F4").Value = 1445
$ws.Range("F5").Value = 1097
$ws.Range("F8").Value = 521
$ws.Range("F9").Value = 222
$ws.Range("F11").Value = 673
$ws.Range("F13").Value = 250
$ws.Range("F17").Value = 151
$ws.Range("F18").Value = 1990
$ws.Range("F20").Value = 430
$ws.Range("F22").Value = 496
$ws.Range("F23").Value = 270
$ws.Range("F29").Value = 282
$ws.Range("F31").Value = 663
$ws.Range("F36").Value = 51
$ws.Range("F37").Value = 241
$ws.Range("F38").Value = 960
$ws.Range("F40").Value = 1561
$ws.Range("F41").Value = 303
